$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 3467225154.74
$ws.Range("P2").Value = 1242555048.92
$ws.Range("Q2").Value = 394748388.05
$ws.Range("R2").Value = 30.8705338516
$ws.Range("S2").Value = 716983053.21
$ws.Range("T2").Value = 11.2072329463
$ws.Range("U2").Value = 556680049.79
$ws.Range("V2").Value = -20.4803009338
$ws.Range("W2").Value = 2342474094.76
$ws.Range("X2").Value = 590706181.27
$ws.Range("Y2").Value = -22.6467093592
$ws.Range("Z2").Value = 41475747.36
$ws.Range("AA2").Value = -23.7585231332
$ws.Range("AB2").Value = 1124751059.98
$ws.Range("AC2").Value = 33.4608276186
$ws.Range("AD2").Value = 1.429294683
$ws.Range("AE2").Value = -9.051648258
$ws.Range("AF2").Value = 85.39045962989999
$ws.Range("AG2").Value = 67.56048396680001
